{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that ends the \"Requisitos\" section\n// (the one containing \"LOQ4233: Gestao de Negocios (Requisito fraco)\").\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"LOQ4233\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The three paragraphs right after the anchor are the trailing blank\n  // spacer paragraph plus the \"Ver no Jupiter...\" and \"(c) 2020 ...\"\n  // footer paragraphs scraped from the site chrome - remove all three,\n  // leaving the anchor directly followed by the page's closing blank\n  // paragraph / page break. Delete from the end backwards so earlier\n  // indices stay valid.\n  paragraphs.items[anchorIndex + 3].delete();\n  paragraphs.items[anchorIndex + 2].delete();\n  paragraphs.items[anchorIndex + 1].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the \"Requisitos\" section\n# (the one containing \"LOQ4233: Gestao de Negocios (Requisito fraco)\").\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*LOQ4233*\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -ne $null) {\n    # The three paragraphs right after the anchor are the trailing\n    # blank spacer paragraph plus the \"Ver no Jupiter...\" and\n    # \"(c) 2020 ...\" footer paragraphs scraped from the site chrome -\n    # remove all three, leaving the anchor directly followed by the\n    # page's closing blank paragraph / page break.\n    $first = $anchor.Next()\n    $last = $first.Next().Next()\n\n    $rng = $d.Range($first.Range.Start, $last.Range.End)\n    $rng.Delete()\n}\n"}
